$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 8 (copy style from A7 so the date format style index is reused)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A8").Value = 44312
$ws.Range("B8").Value = 84
$ws.Range("C8").Value = -1723

# Add row 9
$ws.Range("A9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A9").Value = 44313
$ws.Range("B9").Value = 96
$ws.Range("C9").Value = -1553

# Update selection to match the diff
$ws.Range("C9").Select()
